# Update Leve profit-tracking figures across sheets per scheduled market-price refresh.
$wb = $excel.ActiveWorkbook

$updates = @{
    "ALC" = @{
        51 = @{ "H"=4549.9165; "I"=4400; "J"=4999.6665; "K"=4400; "L"=4999.6665; "M"=-3916 }
        76 = @{ "H"=2714; "I"=2642.5; "J"=3000; "K"=2642.5; "L"=3000; "M"=-2327.5 }
        79 = @{ "H"=2714; "I"=2642.5; "J"=3000; "K"=2642.5; "L"=3000; "M"=-1550.5 }
        80 = @{ "H"=172.16667; "I"=96; "J"=210.25; "K"=288; "L"=630.75; "M"=710 }
        83 = @{ "H"=172.16667; "I"=96; "J"=210.25; "K"=864; "L"=1892.25; "M"=4128 }
        98 = @{ "H"=2155.889; "I"=1694.5294; "J"=9999; "K"=1694.5294; "L"=9999; "M"=-196.5293999999999; "N"=-12995 }
        122 = @{ "H"=2155.889; "I"=1694.5294; "J"=9999; "K"=5083.5882; "L"=29997; "M"=-2633.5882; "N"=-34897 }
    }
    "ARM" = @{
        32 = @{ "H"=3049.65; "I"=1273.9623; "J"=16494.143; "K"=1273.9623; "L"=16494.143; "M"=-986.9622999999999 }
        61 = @{ "H"=50001740; "I"=62501716; "J"=1844; "K"=62501716; "L"=1844; "M"=-62501504; "N"=-2268 }
        88 = @{ "H"=2933.2856; "I"=4999; "J"=2589; "K"=4999; "L"=2589; "M"=-4593; "N"=-3401 }
        91 = @{ "H"=2933.2856; "I"=4999; "J"=2589; "K"=4999; "L"=2589; "M"=-3595; "N"=-5397 }
        132 = @{ "H"=4001779.8; "I"=4546960.5; "J"=3789; "K"=13640881.5; "L"=11367; "M"=-13638351.5 }
        136 = @{ "H"=50001740; "I"=62501716; "J"=1844; "K"=187505148; "L"=5532; "M"=-187502598; "N"=-10632 }
        139 = @{ "H"=50000; "I"=0; "J"=50000; "K"=0; "L"=50000; "N"=-60280 }
    }
    "BSM" = @{
        86 = @{ "H"=3876.5454; "I"=3594.25; "J"=4629.3335; "K"=3594.25; "L"=4629.3335; "M"=-2471.25; "N"=-6875.3335 }
        89 = @{ "H"=3876.5454; "I"=3594.25; "J"=4629.3335; "K"=17971.25; "L"=23146.6675; "M"=-12355.25; "N"=-34378.6675 }
        134 = @{ "H"=42502604; "I"=85001200; "J"=4006.6667; "K"=255003600; "L"=12020.0001; "M"=-255001065; "N"=-17090.0001 }
    }
    "CRP" = @{
        22 = @{ "H"=706.8; "I"=729.7778; "J"=500; "K"=729.7778; "L"=500; "M"=-379.7778 }
        99 = @{ "H"=3164.9092; "I"=2881.4; "J"=6000; "K"=2881.4; "L"=6000; "M"=-1383.4; "N"=-8996 }
        126 = @{ "H"=3164.9092; "I"=2881.4; "J"=6000; "K"=8644.200000000001; "L"=18000; "M"=-6174.200000000001; "N"=-22940 }
    }
    "CUL" = @{
        68 = @{ "H"=2925.5483; "I"=1054.5555; "J"=3243.2642; "K"=3163.6665; "L"=9729.792600000001; "M"=-2352.6665; "N"=-11351.7926 }
        71 = @{ "H"=2925.5483; "I"=1054.5555; "J"=3243.2642; "K"=9490.9995; "L"=29189.3778; "M"=-5434.9995; "N"=-37301.3778 }
    }
    "GSM" = @{
        15 = @{ "H"=59999; "I"=0; "J"=59999; "K"=0; "L"=59999; "N"=-60575 }
        81 = @{ "H"=59999; "I"=0; "J"=59999; "K"=0; "L"=59999; "N"=-61995 }
        84 = @{ "H"=59999; "I"=0; "J"=59999; "K"=0; "L"=179997; "N"=-189981 }
        111 = @{ "H"=100293; "I"=0; "J"=100293; "K"=0; "L"=100293; "N"=-106427 }
        113 = @{ "H"=52475.05; "I"=73776.57000000001; "J"=2771.5; "K"=73776.57000000001; "L"=2771.5; "M"=-71606.57000000001; "N"=-7111.5 }
        132 = @{ "H"=6581388; "I"=7355339; "J"=2807; "K"=22066017; "L"=8421; "M"=-22063487 }
    }
    "LTW" = @{
        22 = @{ "H"=1283.9524; "I"=1468.3889; "J"=177.33333; "K"=1468.3889; "L"=177.33333; "M"=-1173.3889; "N"=-767.3333299999999 }
        27 = @{ "H"=1283.9524; "I"=1468.3889; "J"=177.33333; "K"=1468.3889; "L"=177.33333; "M"=-1361.3889; "N"=-391.33333 }
        40 = @{ "H"=3989.375; "I"=3169.1667; "J"=6450; "K"=3169.1667; "L"=6450; "M"=-3033.1667 }
        132 = @{ "H"=5959446; "I"=9622097; "J"=7638; "K"=28866291; "L"=22914; "M"=-28863761; "N"=-27974 }
    }
    "WVR" = @{
        41 = @{ "H"=34165.668; "I"=28000; "J"=35398.8; "K"=28000; "L"=35398.8; "M"=-27610; "N"=-36178.8 }
        81 = @{ "H"=3566.5; "I"=3566.5; "J"=0; "K"=7133; "L"=0; "M"=-6072 }
        84 = @{ "H"=3566.5; "I"=3566.5; "J"=0; "K"=35665; "L"=0; "M"=-30361 }
        126 = @{ "H"=1723.7778; "I"=1824.1428; "J"=1372.5; "K"=5472.428400000001; "L"=4117.5; "M"=-3002.428400000001; "N"=-9057.5 }
        132 = @{ "H"=20841158; "I"=26318330; "J"=27898; "K"=78954990; "L"=83694; "M"=-78952460 }
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $sheetRows = $updates[$sheetName]
    foreach ($rowNum in $sheetRows.Keys) {
        $rowVals = $sheetRows[$rowNum]
        foreach ($col in $rowVals.Keys) {
            $ws.Range("$col$rowNum").Value = $rowVals[$col]
        }
    }
}

Write-Output "Updated leve profit figures across $($updates.Keys.Count) sheets."